$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1736
$ws.Range("A1736").Value = 44494
$ws.Range("B1736").Value = "'8251330"
$ws.Range("C1736").Value = 3011
$ws.Range("D1736").Value = "Order 8251330 Swish +46727815808"
$ws.Range("F1736").Value = 924.11
$ws.Range("A1736").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1737
$ws.Range("A1737").Value = 44494
$ws.Range("B1737").Value = "'8251330"
$ws.Range("C1737").Value = 2611
$ws.Range("D1737").Value = "Order 8251330 Swish +46727815808"
$ws.Range("F1737").Value = 110.89
$ws.Range("A1737").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1738
$ws.Range("A1738").Value = 44494
$ws.Range("B1738").Value = "'8251330"
$ws.Range("C1738").Value = 1930
$ws.Range("D1738").Value = "Order 8251330 Swish +46727815808"
$ws.Range("E1738").Value = 1035
$ws.Range("A1738").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1739
$ws.Range("A1739").Value = 44494
$ws.Range("C1739").Value = 4010
$ws.Range("D1739").Value = "M&S RB BROMMA K6885"
$ws.Range("E1739").Value = 1023.71
$ws.Range("A1739").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1740
$ws.Range("A1740").Value = 44494
$ws.Range("C1740").Value = 2645
$ws.Range("D1740").Value = "M&S RB BROMMA K6885"
$ws.Range("E1740").Value = 122.85
$ws.Range("A1740").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1741
$ws.Range("A1741").Value = 44494
$ws.Range("C1741").Value = 1930
$ws.Range("D1741").Value = "M&S RB BROMMA K6885"
$ws.Range("F1741").Value = 1146.56
$ws.Range("A1741").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1742
$ws.Range("A1742").Value = 44495
$ws.Range("B1742").Value = "'0261144"
$ws.Range("C1742").Value = 3011
$ws.Range("D1742").Value = "Order 0261144 Swish +46793351577"
$ws.Range("F1742").Value = 739.29
$ws.Range("A1742").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1743
$ws.Range("A1743").Value = 44495
$ws.Range("B1743").Value = "'0261144"
$ws.Range("C1743").Value = 2611
$ws.Range("D1743").Value = "Order 0261144 Swish +46793351577"
$ws.Range("F1743").Value = 88.70999999999999
$ws.Range("A1743").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1744
$ws.Range("A1744").Value = 44495
$ws.Range("B1744").Value = "'0261144"
$ws.Range("C1744").Value = 1930
$ws.Range("D1744").Value = "Order 0261144 Swish +46793351577"
$ws.Range("E1744").Value = 828
$ws.Range("A1744").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1745
$ws.Range("A1745").Value = 44496
$ws.Range("C1745").Value = 4010
$ws.Range("D1745").Value = "Fresh Life"
$ws.Range("E1745").Value = 133.93
$ws.Range("A1745").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1746
$ws.Range("A1746").Value = 44496
$ws.Range("C1746").Value = 2645
$ws.Range("D1746").Value = "Fresh Life"
$ws.Range("E1746").Value = 16.07
$ws.Range("A1746").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1747
$ws.Range("A1747").Value = 44496
$ws.Range("C1747").Value = 1930
$ws.Range("D1747").Value = "Fresh Life"
$ws.Range("F1747").Value = 150
$ws.Range("A1747").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1748
$ws.Range("A1748").Value = 44497
$ws.Range("B1748").Value = "'9281259"
$ws.Range("C1748").Value = 3011
$ws.Range("D1748").Value = "Order 9281259 Card(Stripe)"
$ws.Range("F1748").Value = 883.9299999999999
$ws.Range("A1748").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1749
$ws.Range("A1749").Value = 44497
$ws.Range("B1749").Value = "'9281259"
$ws.Range("C1749").Value = 2611
$ws.Range("D1749").Value = "Order 9281259 Card(Stripe)"
$ws.Range("F1749").Value = 106.07
$ws.Range("A1749").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1750
$ws.Range("A1750").Value = 44497
$ws.Range("B1750").Value = "'9281259"
$ws.Range("C1750").Value = 1930
$ws.Range("D1750").Value = "Order 9281259 Card(Stripe)"
$ws.Range("E1750").Value = 990
$ws.Range("A1750").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1751
$ws.Range("A1751").Value = 44497
$ws.Range("C1751").Value = 4010
$ws.Range("D1751").Value = "WILLYS RISSNE K0135"
$ws.Range("E1751").Value = 180.78
$ws.Range("A1751").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1752
$ws.Range("A1752").Value = 44497
$ws.Range("C1752").Value = 2645
$ws.Range("D1752").Value = "WILLYS RISSNE K0135"
$ws.Range("E1752").Value = 21.69
$ws.Range("A1752").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1753
$ws.Range("A1753").Value = 44497
$ws.Range("C1753").Value = 1930
$ws.Range("D1753").Value = "WILLYS RISSNE K0135"
$ws.Range("F1753").Value = 202.47
$ws.Range("A1753").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1754
$ws.Range("A1754").Value = 44498
$ws.Range("C1754").Value = 7010
$ws.Range("D1754").Value = "NEHA OCT LÖN"
$ws.Range("E1754").Value = 2848
$ws.Range("A1754").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1755
$ws.Range("A1755").Value = 44498
$ws.Range("D1755").Value = "NEHA OCT LÖN"
$ws.Range("E1755").Value = 0
$ws.Range("A1755").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1756
$ws.Range("A1756").Value = 44498
$ws.Range("C1756").Value = 1930
$ws.Range("D1756").Value = "NEHA OCT LÖN"
$ws.Range("F1756").Value = 2848
$ws.Range("A1756").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1757
$ws.Range("A1757").Value = 44498
$ws.Range("B1757").Value = "'8291809"
$ws.Range("C1757").Value = 3011
$ws.Range("D1757").Value = "Order 8291809 Swish +46705093344"
$ws.Range("F1757").Value = 423.21
$ws.Range("A1757").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1758
$ws.Range("A1758").Value = 44498
$ws.Range("B1758").Value = "'8291809"
$ws.Range("C1758").Value = 2611
$ws.Range("D1758").Value = "Order 8291809 Swish +46705093344"
$ws.Range("F1758").Value = 50.79
$ws.Range("A1758").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1759
$ws.Range("A1759").Value = 44498
$ws.Range("B1759").Value = "'8291809"
$ws.Range("C1759").Value = 1930
$ws.Range("D1759").Value = "Order 8291809 Swish +46705093344"
$ws.Range("E1759").Value = 474
$ws.Range("A1759").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1760
$ws.Range("A1760").Value = 44499
$ws.Range("C1760").Value = 7010
$ws.Range("D1760").Value = "Sinthu Lön Octob"
$ws.Range("E1760").Value = 2757
$ws.Range("A1760").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1761
$ws.Range("A1761").Value = 44499
$ws.Range("D1761").Value = "Sinthu Lön Octob"
$ws.Range("E1761").Value = 0
$ws.Range("A1761").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1762
$ws.Range("A1762").Value = 44499
$ws.Range("C1762").Value = 1930
$ws.Range("D1762").Value = "Sinthu Lön Octob"
$ws.Range("F1762").Value = 2757
$ws.Range("A1762").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1763
$ws.Range("A1763").Value = 44499
$ws.Range("C1763").Value = 6540
$ws.Range("D1763").Value = "TWILIO.COM DUBLIN"
$ws.Range("E1763").Value = 174.89
$ws.Range("A1763").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1764
$ws.Range("A1764").Value = 44499
$ws.Range("D1764").Value = "TWILIO.COM DUBLIN"
$ws.Range("E1764").Value = 0
$ws.Range("A1764").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1765
$ws.Range("A1765").Value = 44499
$ws.Range("C1765").Value = 1930
$ws.Range("D1765").Value = "TWILIO.COM DUBLIN"
$ws.Range("F1765").Value = 174.89
$ws.Range("A1765").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1766
$ws.Range("A1766").Value = 44500
$ws.Range("B1766").Value = "'8310827"
$ws.Range("C1766").Value = 3011
$ws.Range("D1766").Value = "Order 8310827 Card(Stripe)"
$ws.Range("F1766").Value = 953.5700000000001
$ws.Range("A1766").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1767
$ws.Range("A1767").Value = 44500
$ws.Range("B1767").Value = "'8310827"
$ws.Range("C1767").Value = 2611
$ws.Range("D1767").Value = "Order 8310827 Card(Stripe)"
$ws.Range("F1767").Value = 114.43
$ws.Range("A1767").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1768
$ws.Range("A1768").Value = 44500
$ws.Range("B1768").Value = "'8310827"
$ws.Range("C1768").Value = 1930
$ws.Range("D1768").Value = "Order 8310827 Card(Stripe)"
$ws.Range("E1768").Value = 1068
$ws.Range("A1768").NumberFormat = "YYYY-MM-DD HH:MM:SS"

